# Resize the scoresheet table: widen the "Div" column (so weight classes
# like "120+" fit) and rebalance the other column widths; set an explicit
# header row height; update the overall table width to match.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- per-column widths (twips -> points; Word COM widths are in points) ---
$t.Columns(1).Width  = 2689 / 20.0   # Name       : 3114 -> 2689
$t.Columns(2).Width  = 850  / 20.0   # Team       : 1276 -> 850
$t.Columns(3).Width  = 3291 / 20.0   # Div        : 1701 -> 3291
$t.Columns(4).Width  = 678  / 20.0   # Lot        : 850  -> 678
$t.Columns(5).Width  = 1105 / 20.0   # Bwt (kg)   : 1001 -> 1105
$t.Columns(6).Width  = 1221 / 20.0   # WtCls (kg) : 1267 -> 1221
$t.Columns(7).Width  = 880  / 20.0   # SQ-1       : 913  -> 880
$t.Columns(8).Width  = 880  / 20.0   # SQ-2       : 913  -> 880
$t.Columns(9).Width  = 881  / 20.0   # SQ-3       : 914  -> 881
$t.Columns(10).Width = 965  / 20.0   # Best SQ    : 1002 -> 965
$t.Columns(11).Width = 965  / 20.0   # BP-1       : 1002 -> 965
$t.Columns(12).Width = 964  / 20.0   # BP-2       : 1001 -> 964
$t.Columns(13).Width = 965  / 20.0   # BP-3       : 1002 -> 965
$t.Columns(14).Width = 965  / 20.0   # Best BP    : 1002 -> 965
$t.Columns(15).Width = 1138 / 20.0   # Sub Total  : 1181 -> 1138
$t.Columns(16).Width = 908  / 20.0   # DL-1       : 942  -> 908
$t.Columns(17).Width = 908  / 20.0   # DL-2       : 942  -> 908
$t.Columns(18).Width = 908  / 20.0   # DL-3       : 942  -> 908
$t.Columns(19).Width = 965  / 20.0   # Best DL    : 1002 -> 965

# --- overall preferred table width: 21967 -> 22126 twips ---
$t.PreferredWidthType = 3
$t.PreferredWidth = 22126 / 20.0

# --- explicit header row height (297 twips, auto rule) ---
$t.Rows(1).Height = 297 / 20.0
